$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Repeatability testing will be done ..." -> split "will be" into its own
#    run and change it to "was" (three runs in the end, matching the diff):
#      "Repeatability testing " | "was" | " done using a laser diode ..."
# ---------------------------------------------------------------------------

# Locate "will be" robustly via Find rather than hard-coded offsets.
$locate = $d.Content
$found = $locate.Find.Execute("testing will be done", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'testing will be done' text"
}

$willBeStart = $locate.Start + 8          # skip "testing "
$willBeEnd = $willBeStart + 7             # "will be" is 7 characters long

# Drop temporary bookmarks right at the two split points. This runtime
# auto-coalesces adjacent, identically-formatted runs whenever a run's text
# is edited, which would otherwise re-merge our new "was" run back into its
# neighbours. Bookmarks act as a hard boundary that blocks that coalescing.
$splitBefore = $d.Range($willBeStart, $willBeStart)
$d.Bookmarks.Add("ZZTmpSplit1", $splitBefore)
$splitAfter = $d.Range($willBeEnd, $willBeEnd)
$d.Bookmarks.Add("ZZTmpSplit2", $splitAfter)

$willBeRange = $d.Range($willBeStart, $willBeEnd)
$willBeRange.Text = "was"

# Remove the temporary barrier bookmarks now that the runs are split.
$d.Bookmarks.Item("ZZTmpSplit1").Delete()
$d.Bookmarks.Item("ZZTmpSplit2").Delete()

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the document (right after
#    the page-break paragraph) to right after the "Altitude only - ..."
#    list item text.
# ---------------------------------------------------------------------------

# Find the "Altitude only ..." bullet text, then expand to the whole
# enclosing paragraph (wdParagraph = 4) so we get real paragraph bounds
# rather than just the matched substring.
$altRange = $d.Content
$altFound = $altRange.Find.Execute("Altitude only", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $altFound) {
    throw "Could not locate 'Altitude only' paragraph"
}
$altRange.Expand(4) | Out-Null
$altParaRange = $altRange.Duplicate
$altParaRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$altParaRange.Collapse(0)                 # collapse to the very end of the text

# A collapsed range that lands exactly on a paragraph's last content
# position (i.e. immediately before the paragraph mark) confuses
# Bookmarks.Add in this runtime (it resolves to document position 0
# instead). Work around it by inserting a unique marker run right after
# the text, locating that marker via Find (so the resulting Range is not
# constructed from the problematic raw offset), adding the bookmark there,
# and then deleting the marker text again.
$altParaRange.InsertAfter("ZZMARKERZZ")

$markerRange = $d.Content
$markerRange.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Collapse(1)                  # collapse to the start of the marker

$d.Bookmarks.Add("_GoBack", $markerRange)

$removeMarker = $d.Content
$removeMarker.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

Write-Output "done"
